$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp update (A1)
$ws.Range("A1").Value = "Datos actualizados a 14 de Julio de 2020 a las 18:54"

# Country rankings shuffled + case counts updated (per 14 Jul 2020 18:54 data refresh)
$updates = @(
    @{ Row=4; Country="Estados Unidos"; Vals=@(3502349, 22866, 1561667, 1802062, 0, 373, 138620) },
    @{ Row=5; Country="Brasil"; Vals=@(1895555, 7596, 1213512, 608882, 0, 240, 73161) },
    @{ Row=6; Country="India"; Vals=@(934565, 26920, 591750, 318514, 0, 574, 24301) },
    @{ Row=11; Country="España"; Vals=@(303699, 666, 0, 0, 0, 3, 28409) },
    @{ Row=12; Country="Reino Unido"; Vals=@(291373, 1240, 0, 0, 0, 138, 44968) },
    @{ Row=16; Country="Italia"; Vals=@(243344, 114, 195441, 12919, 0, 17, 34984) },
    @{ Row=19; Country="Alemania"; Vals=@(200528, 92, 185100, 6289, 0, 0, 9139) },
    @{ Row=30; Country="Suecia"; Vals=@(76001, 34, 0, 0, 0, 11, 5545) },
    @{ Row=56; Country="Irlanda"; Vals=@(25670, 32, 23364, 560, 0, 0, 1746) },
    @{ Row=57; Country="Azerbaiyan"; Vals=@(25113, 543, 16150, 8644, 0, 6, 319) },
    @{ Row=58; Country="Ghana"; Vals=@(24988, 0, 21067, 3782, 0, 0, 139) },
    @{ Row=60; Country="Argelia"; Vals=@(20216, 527, 14295, 4893, 0, 10, 1028) },
    @{ Row=69; Country="Chequia"; Vals=@(13301, 63, 8441, 4505, 0, 2, 355) },
    @{ Row=84; Country="Etiopia"; Vals=@(7969, 203, 2430, 5400, 0, 11, 139) },
    @{ Row=90; Country="Tayikistan"; Vals=@(6643, 47, 5332, 1255, 0, 1, 56) },
    @{ Row=96; Country="Luxemburgo"; Vals=@(5056, 100, 4195, 750, 0, 0, 111) },
    @{ Row=97; Country="Republica de Yibuti"; Vals=@(4979, 2, 4743, 180, 0, 0, 56) },
    @{ Row=100; Country="Grecia"; Vals=@(3883, 57, 1374, 2316, 0, 0, 193) },
    @{ Row=101; Country="Croacia"; Vals=@(3827, 52, 2558, 1149, 0, 1, 120) },
    @{ Row=110; Country="Sri Lanka"; Vals=@(2665, 19, 1988, 666, 0, 0, 11) },
    @{ Row=111; Country="Libano"; Vals=@(2451, 32, 1452, 962, 0, 1, 37) },
    @{ Row=112; Country="Cuba"; Vals=@(2432, 4, 2275, 70, 0, 0, 87) },
    @{ Row=113; Country="Malaui"; Vals=@(2430, 0, 747, 1644, 0, 0, 39) },
    @{ Row=114; Country="Mali"; Vals=@(2423, 11, 1748, 554, 0, 0, 121) },
    @{ Row=125; Country="Sierra Leona"; Vals=@(1651, 9, 1190, 397, 0, 1, 64) },
    @{ Row=133; Country="Tunez"; Vals=@(1306, 4, 1087, 169, 0, 0, 50) },
    @{ Row=136; Country="Jordania"; Vals=@(1198, 15, 1013, 175, 0, 0, 10) },
    @{ Row=156; Country="Reunion"; Vals=@(599, 3, 472, 124, 0, 0, 3) }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 1).Value = $u.Country
    $v = $u.Vals
    for ($c = 0; $c -lt $v.Count; $c++) {
        $ws.Cells.Item($r, 2 + $c).Value = $v[$c]
    }
}
